$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slack bus voltage setpoint (column B) changed from 1.05 to 1.02 pu for all data rows (2-25)
$ws.Range("B2:B25").Value = 1.02

# Recalculated bus voltage magnitudes (columns C-F and I-N) for the 380 kV case
# Row 2
$ws.Cells.Item(2, 3).Value = 1.045513250922339
$ws.Cells.Item(2, 4).Value = 1.040251481188331
$ws.Cells.Item(2, 5).Value = 1.052680320673988
$ws.Cells.Item(2, 6).Value = 1.062165316521385
$ws.Cells.Item(2, 9).Value = 1.043723007886286
$ws.Cells.Item(2, 10).Value = 1.050572583938902
$ws.Cells.Item(2, 11).Value = 1.043034359155837
$ws.Cells.Item(2, 12).Value = 1.055428383217968
$ws.Cells.Item(2, 13).Value = 1.064887403489241
$ws.Cells.Item(2, 14).Value = 1.020656839099321

# Row 3
$ws.Cells.Item(3, 3).Value = 1.046645007966449
$ws.Cells.Item(3, 4).Value = 1.040827656938864
$ws.Cells.Item(3, 5).Value = 1.053694190668258
$ws.Cells.Item(3, 6).Value = 1.063280929839654
$ws.Cells.Item(3, 9).Value = 1.044007418449745
$ws.Cells.Item(3, 10).Value = 1.051351133773937
$ws.Cells.Item(3, 11).Value = 1.043421536742855
$ws.Cells.Item(3, 12).Value = 1.05625463157037
$ws.Cells.Item(3, 13).Value = 1.065817022792727
$ws.Cells.Item(3, 14).Value = 1.02092107119124

# Row 4
$ws.Cells.Item(4, 3).Value = 1.047377291019806
$ws.Cells.Item(4, 4).Value = 1.041200406318878
$ws.Cells.Item(4, 5).Value = 1.054350542399773
$ws.Cells.Item(4, 6).Value = 1.064003232890437
$ws.Cells.Item(4, 9).Value = 1.044190058981392
$ws.Cells.Item(4, 10).Value = 1.051854315161417
$ws.Cells.Item(4, 11).Value = 1.043671290471317
$ws.Cells.Item(4, 12).Value = 1.056788969377763
$ws.Cells.Item(4, 13).Value = 1.066418383711259
$ws.Cells.Item(4, 14).Value = 1.021091707733694

# Row 5
$ws.Cells.Item(5, 3).Value = 1.047685134220924
$ws.Cells.Item(5, 4).Value = 1.041357091188284
$ws.Cells.Item(5, 5).Value = 1.054626547031366
$ws.Cells.Item(5, 6).Value = 1.064306991477714
$ws.Cells.Item(5, 9).Value = 1.044266507616425
$ws.Cells.Item(5, 10).Value = 1.052065711033369
$ws.Cells.Item(5, 11).Value = 1.043776100594398
$ws.Cells.Item(5, 12).Value = 1.057013533325176
$ws.Cells.Item(5, 13).Value = 1.066671156559228
$ws.Cells.Item(5, 14).Value = 1.021163362090687

# Row 6
$ws.Cells.Item(6, 3).Value = 1.047736821969628
$ws.Cells.Item(6, 4).Value = 1.041383398124286
$ws.Cells.Item(6, 5).Value = 1.05467289378821
$ws.Cells.Item(6, 6).Value = 1.06435799989223
$ws.Cells.Item(6, 9).Value = 1.0442793241204
$ws.Cells.Item(6, 10).Value = 1.052101197034991
$ws.Cells.Item(6, 11).Value = 1.043793687731758
$ws.Cells.Item(6, 12).Value = 1.057051234402859
$ws.Cells.Item(6, 13).Value = 1.066713595935048
$ws.Cells.Item(6, 14).Value = 1.021175388407698

# Row 7
$ws.Cells.Item(7, 3).Value = 1.047381404468432
$ws.Cells.Item(7, 4).Value = 1.041202500025643
$ws.Cells.Item(7, 5).Value = 1.0543542300948
$ws.Cells.Item(7, 6).Value = 1.064007291324934
$ws.Cells.Item(7, 9).Value = 1.044191081802216
$ws.Cells.Item(7, 10).Value = 1.051857140399429
$ws.Cells.Item(7, 11).Value = 1.04367269168222
$ws.Cells.Item(7, 12).Value = 1.056791970293592
$ws.Cells.Item(7, 13).Value = 1.066421761428524
$ws.Cells.Item(7, 14).Value = 1.021092665501983

# Row 8
$ws.Cells.Item(8, 3).Value = 1.045895742187765
$ws.Cells.Item(8, 4).Value = 1.04044621735228
$ws.Cells.Item(8, 5).Value = 1.053022898670718
$ws.Cells.Item(8, 6).Value = 1.062542255077288
$ws.Cells.Item(8, 9).Value = 1.043819413915175
$ws.Cells.Item(8, 10).Value = 1.050835821310672
$ws.Cells.Item(8, 11).Value = 1.04316536799061
$ws.Cells.Item(8, 12).Value = 1.055707679971465
$ws.Cells.Item(8, 13).Value = 1.065201606958359
$ws.Cells.Item(8, 14).Value = 1.020746207813086

# Row 9
$ws.Cells.Item(9, 3).Value = 1.043277475217691
$ws.Cells.Item(9, 4).Value = 1.039113023245888
$ws.Cells.Item(9, 5).Value = 1.050679290499605
$ws.Cells.Item(9, 6).Value = 1.059963929902448
$ws.Cells.Item(9, 9).Value = 1.043153832746422
$ws.Cells.Item(9, 10).Value = 1.049031578068876
$ws.Cells.Item(9, 11).Value = 1.042265479270973
$ws.Cells.Item(9, 12).Value = 1.053794714549351
$ws.Cells.Item(9, 13).Value = 1.063050266079249
$ws.Cells.Item(9, 14).Value = 1.020133107162723

# Row 10
$ws.Cells.Item(10, 3).Value = 1.041531679032901
$ws.Cells.Item(10, 4).Value = 1.03822393334716
$ws.Cells.Item(10, 5).Value = 1.049118464666362
$ws.Cells.Item(10, 6).Value = 1.058247216080013
$ws.Cells.Item(10, 9).Value = 1.042702956454404
$ws.Cells.Item(10, 10).Value = 1.047825672698186
$ws.Cells.Item(10, 11).Value = 1.041661604842483
$ws.Cells.Item(10, 12).Value = 1.052517838103132
$ws.Cells.Item(10, 13).Value = 1.061615166727257
$ws.Cells.Item(10, 14).Value = 1.019722626498925

# Row 11
$ws.Cells.Item(11, 3).Value = 1.040775649095648
$ws.Cells.Item(11, 4).Value = 1.037838890373765
$ws.Cells.Item(11, 5).Value = 1.048442980795152
$ws.Cells.Item(11, 6).Value = 1.057504369330996
$ws.Cells.Item(11, 9).Value = 1.042506026119915
$ws.Cells.Item(11, 10).Value = 1.047302767747851
$ws.Cells.Item(11, 11).Value = 1.041399190216527
$ws.Cells.Item(11, 12).Value = 1.051964560898552
$ws.Cells.Item(11, 13).Value = 1.06099354147093
$ws.Cells.Item(11, 14).Value = 1.0195444691228

# Row 12
$ws.Cells.Item(12, 3).Value = 1.040494811336653
$ws.Cells.Item(12, 4).Value = 1.037695860043416
$ws.Cells.Item(12, 5).Value = 1.048192130271864
$ws.Cells.Item(12, 6).Value = 1.057228517773229
$ws.Cells.Item(12, 9).Value = 1.042432622364169
$ws.Cells.Item(12, 10).Value = 1.047108425873236
$ws.Cells.Item(12, 11).Value = 1.041301578015719
$ws.Cells.Item(12, 12).Value = 1.051758991371984
$ws.Cells.Item(12, 13).Value = 1.060762608957348
$ws.Cells.Item(12, 14).Value = 1.01947823088684

# Row 13
$ws.Cells.Item(13, 3).Value = 1.040555052669151
$ws.Cells.Item(13, 4).Value = 1.037726540891619
$ws.Cells.Item(13, 5).Value = 1.048245936112025
$ws.Cells.Item(13, 6).Value = 1.057287685493068
$ws.Cells.Item(13, 9).Value = 1.042448379260493
$ws.Cells.Item(13, 10).Value = 1.047150117927823
$ws.Cells.Item(13, 11).Value = 1.041322522479798
$ws.Cells.Item(13, 12).Value = 1.051803089353356
$ws.Cells.Item(13, 13).Value = 1.060812146276523
$ws.Cells.Item(13, 14).Value = 1.019492442049721

# Row 14
$ws.Cells.Item(14, 3).Value = 1.040752435263535
$ws.Cells.Item(14, 4).Value = 1.03782706760429
$ws.Cells.Item(14, 5).Value = 1.04842224430818
$ws.Cells.Item(14, 6).Value = 1.057481565832577
$ws.Cells.Item(14, 9).Value = 1.042499963745551
$ws.Cells.Item(14, 10).Value = 1.047286705669411
$ws.Cells.Item(14, 11).Value = 1.041391124409312
$ws.Cells.Item(14, 12).Value = 1.051947569637837
$ws.Cells.Item(14, 13).Value = 1.060974453197791
$ws.Cells.Item(14, 14).Value = 1.019538995128104

# Row 15
$ws.Cells.Item(15, 3).Value = 1.040874047225068
$ws.Cells.Item(15, 4).Value = 1.0378890043685
$ws.Cells.Item(15, 5).Value = 1.048530880784242
$ws.Cells.Item(15, 6).Value = 1.057601031803293
$ws.Cells.Item(15, 9).Value = 1.04253171285218
$ws.Cells.Item(15, 10).Value = 1.047370847068833
$ws.Cells.Item(15, 11).Value = 1.041433373820748
$ws.Cells.Item(15, 12).Value = 1.052036581045051
$ws.Cells.Item(15, 13).Value = 1.06107445143248
$ws.Cells.Item(15, 14).Value = 1.019567669711662

# Row 16
$ws.Cells.Item(16, 3).Value = 1.0415818522745
$ws.Cells.Item(16, 4).Value = 1.038249486174792
$ws.Cells.Item(16, 5).Value = 1.049163301939228
$ws.Cells.Item(16, 6).Value = 1.058296526893165
$ws.Cells.Item(16, 9).Value = 1.042715990301629
$ws.Cells.Item(16, 10).Value = 1.047860360542806
$ws.Cells.Item(16, 11).Value = 1.041679000817944
$ws.Cells.Item(16, 12).Value = 1.052554549217096
$ws.Cells.Item(16, 13).Value = 1.061656417317626
$ws.Cells.Item(16, 14).Value = 1.01973444143495

# Row 17
$ws.Cells.Item(17, 3).Value = 1.042025815056229
$ws.Cells.Item(17, 4).Value = 1.038475591192474
$ws.Cells.Item(17, 5).Value = 1.04956010017806
$ws.Cells.Item(17, 6).Value = 1.058732926496984
$ws.Cells.Item(17, 9).Value = 1.042831128118052
$ws.Cells.Item(17, 10).Value = 1.04816722093421
$ws.Cells.Item(17, 11).Value = 1.0418328266687
$ws.Cells.Item(17, 12).Value = 1.052879354567235
$ws.Cells.Item(17, 13).Value = 1.062021410488486
$ws.Cells.Item(17, 14).Value = 1.019838941359744

# Row 18
$ws.Cells.Item(18, 3).Value = 1.042284762564613
$ws.Cells.Item(18, 4).Value = 1.038607468456547
$ws.Cells.Item(18, 5).Value = 1.049791581028332
$ws.Cells.Item(18, 6).Value = 1.058987519444251
$ws.Cells.Item(18, 9).Value = 1.042898122231141
$ws.Cells.Item(18, 10).Value = 1.048346136047852
$ws.Cells.Item(18, 11).Value = 1.041922460632689
$ws.Cells.Item(18, 12).Value = 1.053068771275608
$ws.Cells.Item(18, 13).Value = 1.06223428385205
$ws.Cells.Item(18, 14).Value = 1.019899854212801

# Row 19
$ws.Cells.Item(19, 3).Value = 1.042373055574568
$ws.Cells.Item(19, 4).Value = 1.038652434153929
$ws.Cells.Item(19, 5).Value = 1.049870515991414
$ws.Cells.Item(19, 6).Value = 1.059074337293835
$ws.Cells.Item(19, 9).Value = 1.042920937703463
$ws.Cells.Item(19, 10).Value = 1.048407129389819
$ws.Cells.Item(19, 11).Value = 1.041953008209891
$ws.Cells.Item(19, 12).Value = 1.053133351257391
$ws.Cells.Item(19, 13).Value = 1.062306864624827
$ws.Cells.Item(19, 14).Value = 1.019920617107704

# Row 20
$ws.Cells.Item(20, 3).Value = 1.041978182936427
$ws.Cells.Item(20, 4).Value = 1.038451332874114
$ws.Cells.Item(20, 5).Value = 1.049517523849405
$ws.Cells.Item(20, 6).Value = 1.058686099942934
$ws.Cells.Item(20, 9).Value = 1.042818791865125
$ws.Cells.Item(20, 10).Value = 1.048134305103139
$ws.Cells.Item(20, 11).Value = 1.041816331921926
$ws.Cells.Item(20, 12).Value = 1.0528445098322
$ws.Cells.Item(20, 13).Value = 1.061982252348793
$ws.Cells.Item(20, 14).Value = 1.019827733663806

# Row 21
$ws.Cells.Item(21, 3).Value = 1.040694311421637
$ws.Cells.Item(21, 4).Value = 1.037797465212329
$ws.Cells.Item(21, 5).Value = 1.048370324444306
$ws.Cells.Item(21, 6).Value = 1.057424470844663
$ws.Cells.Item(21, 9).Value = 1.042484780441089
$ws.Cells.Item(21, 10).Value = 1.047246487071364
$ws.Cells.Item(21, 11).Value = 1.041370926704419
$ws.Cells.Item(21, 12).Value = 1.051905025392653
$ws.Cells.Item(21, 13).Value = 1.060926658775791
$ws.Cells.Item(21, 14).Value = 1.019525288136978

# Row 22
$ws.Cells.Item(22, 3).Value = 1.039887005477483
$ws.Cells.Item(22, 4).Value = 1.037386305541115
$ws.Cells.Item(22, 5).Value = 1.047649348379218
$ws.Cells.Item(22, 6).Value = 1.056631667044069
$ws.Cells.Item(22, 9).Value = 1.042273298512085
$ws.Cells.Item(22, 10).Value = 1.046687634152747
$ws.Cells.Item(22, 11).Value = 1.041090074581962
$ws.Cells.Item(22, 12).Value = 1.051313999956198
$ws.Cells.Item(22, 13).Value = 1.060262772777685
$ws.Cells.Item(22, 14).Value = 1.019334766096443

# Row 23
$ws.Cells.Item(23, 3).Value = 1.040314982106711
$ws.Cells.Item(23, 4).Value = 1.037604273208188
$ws.Cells.Item(23, 5).Value = 1.048031521837369
$ws.Cells.Item(23, 6).Value = 1.057051906534858
$ws.Cells.Item(23, 9).Value = 1.042385548951706
$ws.Cells.Item(23, 10).Value = 1.046983954112745
$ws.Cells.Item(23, 11).Value = 1.041239036063267
$ws.Cells.Item(23, 12).Value = 1.051627345561148
$ws.Cells.Item(23, 13).Value = 1.060614729712724
$ws.Cells.Item(23, 14).Value = 1.019435799823756

# Row 24
$ws.Cells.Item(24, 3).Value = 1.041999705866848
$ws.Cells.Item(24, 4).Value = 1.0384622941825
$ws.Cells.Item(24, 5).Value = 1.049536762150818
$ws.Cells.Item(24, 6).Value = 1.05870725869715
$ws.Cells.Item(24, 9).Value = 1.042824366592943
$ws.Cells.Item(24, 10).Value = 1.048149178570621
$ws.Cells.Item(24, 11).Value = 1.041823785466228
$ws.Cells.Item(24, 12).Value = 1.052860254781447
$ws.Cells.Item(24, 13).Value = 1.061999946291162
$ws.Cells.Item(24, 14).Value = 1.019832798063275

# Row 25
$ws.Cells.Item(25, 3).Value = 1.043954406086082
$ws.Cells.Item(25, 4).Value = 1.039457742090566
$ws.Cells.Item(25, 5).Value = 1.051284890051256
$ws.Cells.Item(25, 6).Value = 1.060630104919638
$ws.Cells.Item(25, 9).Value = 1.043327163236728
$ws.Cells.Item(25, 10).Value = 1.049498559102963
$ws.Cells.Item(25, 11).Value = 1.04249882032662
$ws.Cells.Item(25, 12).Value = 1.054289536917761
$ws.Cells.Item(25, 13).Value = 1.06360659177364
$ws.Cells.Item(25, 14).Value = 1.020291916261416
